# register competition create excel done
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "sheet1"

# 2. Populate the hidden category list in column Z (used as the source range
#    for the dropdown data validations below).
$categories = @(
    "זכר דואן בנים 7-7 (קוד: 16)",
    "זכר נשק קצר חלקי 8-9 (קוד: 11)",
    "זכר צאן סואן בנים 9-12 (קוד: 15)",
    "זכר צאן חלקי בנים 9-11 (קוד: 1)",
    "זכר נשק ארוך בנים 10-12 (קוד: 3)",
    "זכר נשק קצר בנים 11-11 (קוד: 8)",
    "זכר נשק ארוך בנים 14-14 (קוד: 5)",
    "זכר צאן חלקי בנים 14-14 (קוד: 2)",
    "מעורב נשק ארוך 14-14 (קוד: 6)",
    "נקבה דואילין בנות  (קוד: 18)",
    "נקבה דואן בנות 7-7 (קוד: 17)",
    "נקבה בלי נשק בנות 9-12 (קוד: 12)",
    "נקבה נשק קצר בנות 9-11 (קוד: 9)",
    "נקבה נשק ארוך בנות 10-10 (קוד: 4)",
    "נקבה נשק ארוך חלקי 11-11 (קוד: 7)",
    "נקבה נשק קצר בנות 14-14 (קוד: 10)",
    "נקבה צאן סואן בנות 14-14 (קוד: 14)",
    "נקבה בלי נשק בוגרות 18+ (קוד: 13)"
)

# First cell gets the formatting directly (white text, right aligned, 10pt)
# so it creates exactly one new font + one new cell style ...
$ws.Range("Z1").Value = $categories[0]
$ws.Range("Z1").Font.Size = 10
$ws.Range("Z1").Font.Color = 16777215
$ws.Range("Z1").Font.Name = "Calibri"
$ws.Range("Z1").HorizontalAlignment = -4152

# ... remaining cells just get their value ...
for ($i = 1; $i -lt $categories.Length; $i++) {
    $ws.Cells.Item($i + 1, 26).Value = $categories[$i]
}

# ... then the format from Z1 is copied onto Z2:Z18 so every cell reuses the
# same style index instead of minting a new one each time.
$ws.Range("Z1").Copy()
$ws.Range("Z2:Z18").PasteSpecial(-4122)

# 3. Point the three dropdown validations at the Z column list instead of the
#    inline comma separated literal, and require a selection (allowBlank=0).
$ranges = @("F2:F14", "G2:G14", "H2:H14")
foreach ($addr in $ranges) {
    $rng = $ws.Range($addr)
    $rng.Validation.Delete()
    $rng.Validation.Add(3, 1, 1, "=sheet1!`$Z`$1:`$Z`$100")
    $rng.Validation.IgnoreBlank = $false
    $rng.Validation.InCellDropdown = $true
    $rng.Validation.ShowInput = $true
    $rng.Validation.ShowError = $true
    $rng.Validation.InputMessage = "בחר קטגוריה"
    $rng.Validation.ErrorMessage = "Invalid choice was chosen"
}
